$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update password values for a few rows
$ws.Range("B3").Value = 12345
$ws.Range("B4").Value = 1378
$ws.Range("B5").Value = 12678

# Replace the password value on the last row with a non-numeric text value.
# A leading apostrophe forces it to be stored as explicit text (quote-prefixed),
# matching the original cell's text formatting.
$ws.Range("B8").Value = "'wqerewr"

# Move the active selection from D8 to G4
$ws.Range("G4").Select()
